$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.245.81'
$ws.Range('E2').Value = '  +4.44%  '
$ws.Range('D3').Value = '2.510.08'
$ws.Range('E3').Value = '  +2.16%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '495.40'
$ws.Range('E5').Value = '  +2.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.05'
$ws.Range('E6').Value = '  +11.24%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.516'
$ws.Range('E8').Value = '  +3.02%  '
$ws.Range('D9').Value = '2.529.62'
$ws.Range('E10').Value = '  +4.63%  '
$ws.Range('E11').Value = '  +5.51%  '
$ws.Range('E12').Value = '  +4.42%  '
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').Value = '2.946.19'
$ws.Range('E14').Value = '  +3.26%  '
$ws.Range('D15').Value = '57.365.93'
$ws.Range('E15').Value = '  +4.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.43'
$ws.Range('E16').Value = '  +4.76%  '
$ws.Range('E17').Value = '  +3.20%  '
$ws.Range('D18').Value = '2.528.54'
$ws.Range('E18').Value = '  +3.49%  '
$ws.Range('E19').Value = '  +6.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.36'
$ws.Range('E20').Value = '  +4.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.21'
$ws.Range('E21').Value = '  +3.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.94'
$ws.Range('E23').Value = '  +5.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.62'
$ws.Range('E24').Value = '  +2.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.412'
$ws.Range('E25').Value = '  +1.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.165'
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('D28').Value = '2.618.14'
$ws.Range('E28').Value = '  +3.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.65'
$ws.Range('E29').Value = '  +4.22%  '
$ws.Range('D30').Value = '0.0₃0834'
$ws.Range('E30').Value = '  +7.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.71'
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.54'
$ws.Range('E33').Value = '  +4.59%  '
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.32'
$ws.Range('E35').Value = '  +3.70%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.17'
$ws.Range('E36').Value = '  +4.80%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.84'
$ws.Range('E37').Value = '  +6.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.892'
$ws.Range('E38').Value = '  +5.31%  '
$ws.Range('E39').Value = '  +10.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '34.41'
$ws.Range('E40').Value = '  +3.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.56'
$ws.Range('E41').Value = '  +4.61%  '
$ws.Range('E42').Value = '  +4.41%  '
$ws.Range('E43').Value = '  +3.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.995'
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.96'
$ws.Range('E45').Value = '  +6.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '270.06'
$ws.Range('E46').Value = '  +6.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0943'
$ws.Range('E47').Value = '  +4.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0231'
$ws.Range('E48').Value = '  +4.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.22'
$ws.Range('E49').Value = '  +1.32%  '
$ws.Range('E50').Value = '  +6.38%  '
$ws.Range('D51').Value = '1.904.36'
$ws.Range('E51').Value = '  -1.44%  '
